$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing "Total Returned to the White House" row (row 40).
# The sheet's summary block is being restructured from 7 rows (34-40) down to
# 6 rows (34-39), so the last row is dropped entirely.
$ws.Rows(40).Delete()

# --- Relabel the per-category detail rows so each sub-label is prefixed with
# its section name (e.g. "     New nominations" -> "     Civilian, New nominations") ---

# Civilian (section header row 6)
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

# Other Civilian (section header row 12)
$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Unconfirmed "

# Air Force (section header row 16)
$ws.Range("A17").Value = "     Air Force, New nominations"
$ws.Range("A18").Value = "     Air Force, Confirmed "
$ws.Range("A19").Value = "     Air Force, Unconfirmed "
$ws.Range("A20").Value = "     Air Force, Returned to White House "

# Army (section header row 21)
$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Confirmed "
$ws.Range("A24").Value = "     Army, Unconfirmed "
$ws.Range("A25").Value = "     Army, Returned to White House "

# Navy (section header row 26)
$ws.Range("A27").Value = "     Navy, New nominations"
$ws.Range("A28").Value = "     Navy, Confirmed "
$ws.Range("A29").Value = "     Navy, Unconfirmed "

# Marine Corps (section header row 30)
$ws.Range("A31").Value = "     Marine Corps, New nominations"
$ws.Range("A32").Value = "     Marine Corps, Confirmed "
$ws.Range("A33").Value = "     Marine Corps, Unconfirmed "

# --- Rebuild the summary block (rows 34-39) ---
# Row 34 used to just be the "Summary" section header with no value; it now
# becomes a data row carrying the total new-nominations figure that used to
# sit on row 36.
$ws.Range("A34").Value = "Total new nominations"
$ws.Range("B34").Value = 28423
$ws.Range("B34").NumberFormat = "#,##0"

$ws.Range("A35").Value = "Total carryover nominations"
# B35 value (0) is unchanged

$ws.Range("A36").Value = "Total confirmed "
# B36 value (21580) is unchanged

$ws.Range("A37").Value = "Total unconfirmed "
# B37 value (6812) is unchanged

$ws.Range("A38").Value = "Total withdrawn "
$ws.Range("B38").Value = 13
$ws.Range("B38").NumberFormat = "General"

$ws.Range("A39").Value = "Total returned"
# B39 value (18) is unchanged
